$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change E3 from the text "b" to the number 3 (float -> int in middle of object dtype)
$ws.Range("E3").Value = 3

# Update the active selection to E4, matching the author's interaction
$ws.Range("E4").Select()
